# Updated symbol list on Sun Dec 18 06:50:20 UTC 2022 with GitHub Actions
#
# Applies the per-coin price/volume refresh captured in the diff. Prices in
# column D are stored as text (e.g. "0.8030" must stay "0.8030", not become
# the number 0.803), so those cells get NumberFormat "@" (Text) before the
# value is written so trailing zeros / exact formatting survive. Columns
# B/C/E are plain text already and don't need the NumberFormat nudge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# --- Column D (Price) updates -------------------------------------------------
Set-TextValue "D2"  "246.21"
Set-TextValue "D4"  "5.575"
Set-TextValue "D5"  "0.05621"
Set-TextValue "D6"  "3.405"
Set-TextValue "D7"  "6.473"
Set-TextValue "D8"  "0.8031"
Set-TextValue "D9"  "1.070"
Set-TextValue "D10" "0.1428"
Set-TextValue "D11" "0.07403"
Set-TextValue "D12" "0.03186"
Set-TextValue "D13" "0.02964"
Set-TextValue "D14" "0.09254"
Set-TextValue "D15" "0.001670"
Set-TextValue "D16" "3.270"
Set-TextValue "D17" "0.04708"
Set-TextValue "D18" "0.0005758"
Set-TextValue "D20" "0.001051"
Set-TextValue "D23" "0.0004606"
Set-TextValue "D24" "3.982"
Set-TextValue "D25" "2.093"
Set-TextValue "D27" "0.1299"
Set-TextValue "D40" "0.04179"
Set-TextValue "D41" "0.007069"
Set-TextValue "D44" "0.009811"
Set-TextValue "D45" "0.00005626"
Set-TextValue "D47" "0.6809"
Set-TextValue "D48" "0.02853"

# --- Column E (Volume(1h)) text tweaks ----------------------------------------
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# --- Rows 42/43 swap: CEJI <-> BKEXToken rankings shifted ---------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1042"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002974"
$ws.Range("E43").Value = "42CEJICEJI"
